$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

# --- Update the "time_taken" (column F) timestamps on the existing "data" sheet ---
$dataWs.Range("F2").Value = "2021-10-05 14:33:15.574539"
$dataWs.Range("F3").Value = "2021-10-05 14:33:15.574546"
$dataWs.Range("F4").Value = "2021-10-05 14:33:15.574549"
$dataWs.Range("F5").Value = "2021-10-05 14:33:15.574551"
$dataWs.Range("F6").Value = "2021-10-05 14:33:15.574554"
$dataWs.Range("F7").Value = "2021-10-05 14:33:15.574556"
$dataWs.Range("F8").Value = "2021-10-05 14:33:15.574558"
$dataWs.Range("F9").Value = "2021-10-05 14:33:15.574560"
$dataWs.Range("F10").Value = "2021-10-05 14:33:15.574563"
$dataWs.Range("F11").Value = "2021-10-05 14:33:15.574565"
$dataWs.Range("F12").Value = "2021-10-05 14:33:15.574568"
$dataWs.Range("F13").Value = "2021-10-05 14:33:15.574570"
$dataWs.Range("F14").Value = "2021-10-05 14:33:15.574573"
$dataWs.Range("F15").Value = "2021-10-05 14:33:15.574575"
$dataWs.Range("F16").Value = "2021-10-05 14:33:15.574577"
$dataWs.Range("F17").Value = "2021-10-05 14:33:15.574580"
$dataWs.Range("F18").Value = "2021-10-05 14:33:15.574582"
$dataWs.Range("F19").Value = "2021-10-05 14:33:15.574585"
$dataWs.Range("F20").Value = "2021-10-05 14:33:15.574587"

# --- Add a new "metadata" worksheet, positioned right after "data" ---
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataWs)
$newWs.Name = "metadata"

# Reuse the same header-row style (bold, centered, thin border) used on "data"!B1:F1
$dataWs.Range("B1:F1").Copy()
$newWs.Range("B1:G1").PasteSpecial(-4122)

# Reuse the same index-column style used on "data"!A2
$dataWs.Range("A2").Copy()
$newWs.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Header row
$newWs.Range("B1").Value = "data_name"
$newWs.Range("C1").Value = "data_id"
$newWs.Range("D1").Value = "data_version"
$newWs.Range("E1").Value = "data_version_created"
$newWs.Range("F1").Value = "panel_query_time"
$newWs.Range("G1").Value = "panel_get_request"

# Data row
$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "Autonomic neuropathy"
$newWs.Range("C2").Value = 3439

# "0.48" must be stored as text, not as a number - force text then drop the
# number-format styling so the cell ends up unstyled, matching the source data.
$newWs.Range("D2").NumberFormat = "@"
$newWs.Range("D2").Value = "0.48"
$newWs.Range("D2").ClearFormats()

$newWs.Range("E2").Value = "2021-05-21T10:59:00.640114Z"
$newWs.Range("F2").Value = "2021-10-05 14:33:15.571614"
$newWs.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3439/?format=json"
